# Applies the "Trade closed / new trade opened" update to the live trading
# results workbook:
#   - Summary sheet roll-up metrics
#   - Strategy Status roll-up row for MarketMaking
#   - All Trades: trade #127 flips from OPEN -> CLOSED, and a brand new
#     trade #160 is appended
#   - MarketMaking: same two changes, mirrored with that sheet's own
#     column layout

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.4
$summary.Range("B4").Value = 1.19
$summary.Range("B6").Value = 127
$summary.Range("B8").Value = 48
$summary.Range("B9").Value = 43.31

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.4
$status.Range("D5").Value = 94
$status.Range("E5").Value = 1.08
$status.Range("F5").Value = 1.4
$status.Range("G5").Value = 43.62

# ---------------------------------------------------------------------
# All Trades sheet - trade #127 (row 128) closes, trade #160 (row 161) opens
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G128").Value = 0.1
$allTrades.Range("H128").Value = "CLOSED"
$allTrades.Range("I128").Value = -9.0909
$allTrades.Range("J128").Value = -0.01
$allTrades.Range("K128").Value = 101.4
$allTrades.Range("L128").Value = "early_exit"
$allTrades.Range("M128").Value = 0.13

$allTrades.Range("A161").Value = 160
$allTrades.Range("B161").NumberFormat = "@"
$allTrades.Range("B161").Value = "2026-02-17"
$allTrades.Range("C161").Value = "21:27:19"
$allTrades.Range("D161").Value = "MarketMaking"
$allTrades.Range("E161").Value = "UP"
$allTrades.Range("F161").Value = 0.11
$allTrades.Range("H161").Value = "OPEN"
$allTrades.Range("I161").Value = 0
$allTrades.Range("J161").Value = 0
$allTrades.Range("K161").Value = 101.411797784678
$allTrades.Range("M161").Value = 0
$allTrades.Range("N161").Value = 0
$allTrades.Range("O161").Value = 0
$allTrades.Range("P161").Value = 0.6
$allTrades.Range("Q161").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# MarketMaking sheet - trade #127 (row 95) closes, trade #160 (row 128) opens
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Range("G95").Value = 0.1
$mm.Range("H95").Value = "CLOSED"
$mm.Range("I95").Value = -9.0909
$mm.Range("J95").Value = -0.01
$mm.Range("K95").Value = 101.4
$mm.Range("P95").Value = "early_exit"
$mm.Range("Q95").Value = 0.13

$mm.Range("A128").Value = 160
$mm.Range("B128").NumberFormat = "@"
$mm.Range("B128").Value = "2026-02-17"
$mm.Range("C128").Value = "21:27:19"
$mm.Range("D128").Value = "MarketMaking"
$mm.Range("E128").Value = "UP"
$mm.Range("F128").Value = 0.11
$mm.Range("H128").Value = "OPEN"
$mm.Range("I128").Value = 0
$mm.Range("J128").Value = 0
$mm.Range("K128").Value = 101.411797784678
$mm.Range("L128").Value = 0
$mm.Range("M128").Value = 0
$mm.Range("N128").Value = 0.6
$mm.Range("O128").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q128").Value = 0
